# Applies the "Problem/Solution/Conclusions" bold-label formatting pass
# (plus a few wording tweaks) to the Question slides, per the commit diff.

# Rewrites $para's whole text to $fullText (character-exact, incl. any
# smart punctuation) and makes the leading $label bold, leaving the rest
# of the paragraph in regular (non-bold) formatting.
function Set-LabelBold($para, $label, $fullText) {
    $allChars = $para.Characters(1, $para.Text.Length)
    $allChars.Text = $fullText
    $para.Characters(1, $label.Length).Font.Bold = $true
}

# Same as Set-LabelBold, but also underlines the first occurrence of
# $underlineWord inside the paragraph (used for "MapReduce").
function Set-LabelBoldWithUnderline($para, $label, $fullText, $underlineWord) {
    Set-LabelBold $para $label $fullText
    $uStart = $fullText.IndexOf($underlineWord) + 1
    $para.Characters($uStart, $underlineWord.Length).Font.Underline = $true
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3 (Question 1) - Content Placeholder 2
# ---------------------------------------------------------------------
$shp = $p.Slides.Item(3).Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

Set-LabelBold $tr.Paragraphs(1,1) "Problem: " `
    "Problem: Identify the countries where percentage of female graduates is less than 30%."

Set-LabelBoldWithUnderline $tr.Paragraphs(2,1) "Solution: " `
    "Solution: Using MapReduce, find the most recent years for non-null values for gross graduation ratio for females in tertiary education up to five years combined. Take these percentages and find the average by dividing the sum of the percentages by the number of percentages (years accounted for) isolating below 30%." `
    "MapReduce"

Set-LabelBold $tr.Paragraphs(3,1) "Conclusions: " `
    "Conclusions: Out of 89 nations, the vast majority in the final output are Third World. Sub-Saharan Africa are on average on the lower end of the percentage range while Asian and some South American countries are closer to the cap of 30%."

# ---------------------------------------------------------------------
# Slide 4 (Question 3) - Content Placeholder 2
# ---------------------------------------------------------------------
$shp = $p.Slides.Item(4).Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

Set-LabelBold $tr.Paragraphs(1,1) "Problem: " `
    "Problem: List the percentage of change made in male employment from the year 2000."

Set-LabelBoldWithUnderline $tr.Paragraphs(2,1) "Solution: " `
    "Solution: Using MapReduce, find the percent of employment for the year 2000 and compare it to 2016 or the most recent year. Take both of these percentages." `
    "MapReduce"

Set-LabelBold $tr.Paragraphs(3,1) "Conclusions: " `
    "Conclusions: Mostly European countries fell into the negative values for change in male employment since the year 2000 while South American countries’ numbers are rising."

# ---------------------------------------------------------------------
# Slide 5 (Question 2) - Content Placeholder 2
# ---------------------------------------------------------------------
$shp = $p.Slides.Item(5).Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

Set-LabelBold $tr.Paragraphs(1,1) "Problem: " `
    "Problem: List the average increase in female education in the U.S. from the year 2000."

Set-LabelBold $tr.Paragraphs(2,1) "Solution: " `
    "Solution: Find the percent change between the years 2004 and 2015 (2000 – 2003, 2016 null) for women in the US based on master’s, bachelor’s , secondary school and no schooling and average them."

Set-LabelBold $tr.Paragraphs(3,1) "Conclusions: " `
    "Conclusions: Female education in the US dropped overall by 22% between the years 2004 and 2015."

Set-LabelBold $tr.Paragraphs(4,1) "Stipulations: " `
    "Stipulations: Starting from 2004, all levels of tertiary education were included; however 2012 onward, no other notes in the data included."

# ---------------------------------------------------------------------
# Slide 6 (Question 4) - Content Placeholder 2
# ---------------------------------------------------------------------
$shp = $p.Slides.Item(6).Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

Set-LabelBold $tr.Paragraphs(1,1) "Problem: " `
    "Problem: List the percentage of change in female employment from the year 2000."

Set-LabelBold $tr.Paragraphs(2,1) "Solution: " `
    "Solution: Find the percent change in the employment to population ratio for females globally between the years of 2000 and 2016, since these were the most recent years after 1999 available for the “WLD” data set."

Set-LabelBold $tr.Paragraphs(3,1) "Conclusions: " `
    "Conclusions: The percentage of change for the world’s female employment went down about 4%."

# ---------------------------------------------------------------------
# Slide 7 (Question 5) - Content Placeholder 2 (plain wording changes,
# no bold-label restructuring here)
# ---------------------------------------------------------------------
$shp = $p.Slides.Item(7).Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$tr.Paragraphs(1,1).Text = "Problem: List the percentage of women who are financially independent in all countries."
$tr.Paragraphs(2,1).Text = "Solution: "
